$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fix in C3: insert an almanac tooltip right after "is missing." ---
$c3Old = $ws.Range("C3").Value()
$c3New = $c3Old.Replace("is missing. {character:taken:possPronoun:cap}", "is missing.{almanac:[character-taken-name]:your friend from [location-hometown-name], taken to [location-baronhome-name] by the [character-baron-baron]'s men} {character:taken:possPronoun:cap}")
$ws.Range("C3").Value = $c3New

# --- Fix RFW1.3 choice / outro text (rows 35) ---
$ws.Range("E35").Value = "Turn and fight"
$ws.Range("F35").Value = "You stop in your tracks, put your hands up, and slowly turn around. The {character:baron:baron}'s men close the gap between you, smiling, slowing to a walk as they approach.{RFW1.3a}"
$ws.Range("G35").Value = "You turn to face them, {if:item:sword:drawing your sword}{if:noitem:sword:readying yourself for a fight}. The Baron's men just laugh, drawing their swords and advancing quickly to close the distance. They'll be on you in seconds.{RFW1.3a}"

# --- Row height adjustments to fit the new text ---
$ws.Rows.Item(3).RowHeight = 214.15
$ws.Rows.Item(35).RowHeight = 68.65

# --- Update the view / selection (scrolled down & selecting A36) ---
$ws.Range("A36").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
